$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header text for column D (shared string "Antal tasks" -> "Antal tasks tilbage")
$ws.Range("D1").Value = "Antal tasks tilbage"

# Update existing data values
$ws.Range("B2").Value = 10
$ws.Range("D2").Value = 10
$ws.Range("D3").Value = 8
$ws.Range("D4").Value = 6

# Add new rows of data
$ws.Range("C5").Value = 6
$ws.Range("D5").Value = 6
$ws.Range("C6").Value = 6
$ws.Range("D6").Value = 4

# Autofit column D to match new header width
$ws.Columns("D").AutoFit()

# Move selection to D7 like in the target file
$ws.Range("D7").Select()
